$wb = $excel.ActiveWorkbook

# start_price sheet
$ws = $wb.Worksheets.Item("start_price")
$ws.Range("A2").Value = 1751

# Linear sheet
$ws = $wb.Worksheets.Item("Linear")
$ws.Range("B2").Value = -0.07776438606882616
$ws.Range("B3").Value = 0.01429109734191718
$ws.Range("B4").Value = 672.8525460485779
$ws.Range("B5").Value = "[1.0, 0.15021283303220878, -0.044774089363924324, 0.004133520297264346, -0.029385722600000934, -0.04162578194294981, 0.13414897773629786, 0.2798238268943356, 0.14356536574540466, -0.06548117507535693, -0.02232570854942809, -0.029808798604038036, -0.05742458783116447, 0.12036293402667449, 0.2784813569886377, 0.08549830868653784, -0.05578375981896078, -0.05009904974230347, -0.036525588859309774, -0.06925130190606413]"

# NonLinear sheet
$ws = $wb.Worksheets.Item("NonLinear")
$ws.Range("B3").Value = 0.9103699843668578
$ws.Range("B4").Value = -1.470588062475888
$ws.Range("B5").Value = -0.06205060960652477
$ws.Range("B6").Value = 681.2677249578644
$ws.Range("B7").Value = 1.375371575639423
$ws.Range("B8").Value = -0.08574817672285553
$ws.Range("B9").Value = 664.0970834904706
$ws.Range("B10").Value = "[1.0, 0.1446862337135457, -0.040765671922998147, 0.00728610145484793, -0.02568754351490599, -0.039418961063793743, 0.13216065103263858, 0.2703274004515764, 0.14187870835122554, -0.06025347195698553, -0.017855259726688735, -0.024896648033404, -0.05138063414805476, 0.11828746832514847, 0.26997067575452904, 0.08107449525087163, -0.0510528589480201, -0.047995137084212, -0.033935097459040345, -0.06362131275203785]"
